$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add Arabic ("ara") gender rows (8-10) at the bottom of the table, mirroring
# the structure already used for the "fra" rows (5-7): lang_code / code /
# name / is_active.

# --- Row 8 : Male -----------------------------------------------------
$ws.Range("A8").Value = "ara"
$ws.Range("B8").Value = "MLE"
$ws.Range("C8").Value = "ذكر"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4163) | Out-Null   # xlPasteValues -> keep "TRUE" as text, not boolean

# --- Row 9 : Female -----------------------------------------------------
$ws.Range("A9").Value = "ara"
$ws.Range("B9").Value = "FLE"
$ws.Range("C9").Value = "أنثى"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4163) | Out-Null

# --- Row 10 : Others -----------------------------------------------------
$ws.Range("A10").Value = "ara"
$ws.Range("B10").Value = "OTH"
$ws.Range("C10").Value = "آحرون"
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163) | Out-Null  # D4 holds "FALSE"

# The Arabic text needs wrapping, so the "name" column for these new rows
# gets its own style (wrapText = true) and a taller row height.
$wrapRange = $ws.Range("C8:C10")
$wrapRange.WrapText = $true

$ws.Rows.Item(8).RowHeight = 16.4
$ws.Rows.Item(9).RowHeight = 16.4
$ws.Rows.Item(10).RowHeight = 16.4

$ws.Range("C10").Select() | Out-Null
